$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 111111111
$ws.Range("G3").Value = 222222222
$ws.Range("G4").Value = 222222222

$ws.Range("F9").Select()
